$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.917.74'
$ws.Range("E2").Value = '  +2.83%  '
$ws.Range("D3").Value = '3.981.27'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '612.51'
$ws.Range("E5").Value = '  +14.24%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.88'
$ws.Range("E6").Value = '  +11.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.683'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.754'
$ws.Range("E9").Value = '  +2.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  +1.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.41'
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("E12").Value = '  +1.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.98'
$ws.Range("E13").Value = '  +3.73%  '
$ws.Range("D14").Value = '4.627.50'
$ws.Range("E14").Value = '  +1.21%  '
$ws.Range("D15").Value = '3.990.69'
$ws.Range("E15").Value = '  +1.09%  '
$ws.Range("E16").Value = '  +8.93%  '
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.54'
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("E19").Value = '  +0.44%  '
$ws.Range("D20").Value = '72.639.22'
$ws.Range("E20").Value = '  +2.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '439.70'
$ws.Range("E21").Value = '  +4.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.87'
$ws.Range("E22").Value = '  +14.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '96.69'
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("E24").Value = '  -3.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.41'
$ws.Range("E25").Value = '  -0.87%  '
$ws.Range("E26").Value = '  +12.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.44'
$ws.Range("E27").Value = '  +1.33%  '
$ws.Range("E28").Value = '  +1.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.51'
$ws.Range("E29").Value = '  -2.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.48'
$ws.Range("E30").Value = '  +0.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.83'
$ws.Range("E31").Value = '  +1.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.90'
$ws.Range("E32").Value = '  +4.36%  '
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '48.77'
$ws.Range("E34").Value = '  -4.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '71.44'
$ws.Range("E35").Value = '  +8.96%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '656.44'
$ws.Range("E36").Value = '  -3.07%  '
$ws.Range("D37").Value = '0.0₃0899'
$ws.Range("E37").Value = '  +11.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.439'
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.147'
$ws.Range("E39").Value = '  -0.67%  '
$ws.Range("B40").Value = 'ThetaToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.37'
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.34'
$ws.Range("E42").Value = '  +4.95%  '
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("E44").Value = '  +1.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.66'
$ws.Range("E45").Value = '  +5.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.150'
$ws.Range("E46").Value = '  +0.91%  '
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("E48").Value = '  +1.91%  '
$ws.Range("D49").Value = '2.904.02'
$ws.Range("E49").Value = '  +12.23%  '
$ws.Range("E50").Value = '  +2.42%  '
